# Populating the impostor scores for "User 10" sheet, and moving the
# active-tab/tab-selected marker from "User 10" to "User 2".

$wb = $excel.ActiveWorkbook

$ws10 = $wb.Worksheets.Item("User 10")

# Fill in rows 14-23 on sheet "User 10" with the impostor-score formulas.
# Row (13 + i) references row i of the raw data (i = 1..10), comparing it
# against the column averages in row 12.
for ($i = 1; $i -le 10; $i++) {
    $destRow = 13 + $i

    # First cell (column A) gets its own (non-shared) formula, matching the
    # pattern already used for the "Row 12" averages and for "User 1"'s
    # existing impostor-score rows.
    $formulaA = "=ABS(A" + $i + "-A12)/10"
    $ws10.Range("A$destRow").Formula = $formulaA

    # Columns B:G share one formula definition (relative references let
    # Excel expand it per-column), exactly like the B12:G12 shared formula.
    $formulaBG = "=ABS(B" + $i + "-B12)/10"
    $ws10.Range("B$destRow" + ":G$destRow").Formula = $formulaBG
}

# Update selection on "User 10" to reflect the new active cell.
$ws10.Range("K23").Select()

# Move the selected tab from "User 10" to "User 2".
$ws2 = $wb.Worksheets.Item("User 2")
$ws2.Activate()
